$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "1.000", "26.482.61") that must be
# written as TEXT, not auto-converted to a number/date by Excel. $q is a single
# apostrophe char cast to [string] (NOT left as [char], which PowerShell would
# numerically add to a parseable numeric string instead of concatenating).
# Prepending it forces a text entry, same as typing it by hand in the grid.
$q = [string][char]39

$ws.Range("D2").Value = $q + '26.482.61'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = $q + '1.840.19'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("D4").Value = $q + '1.000'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = $q + '260.67'
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("D6").Value = $q + '1.000'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = $q + '0.5370'
$ws.Range("E7").Value = '  +2.28%  '
$ws.Range("D8").Value = $q + '0.2949'
$ws.Range("E8").Value = '  -8.86%  '
$ws.Range("D9").Value = $q + '0.06953'
$ws.Range("E9").Value = '  +2.20%  '
$ws.Range("D10").Value = $q + '17.31'
$ws.Range("E10").Value = '  -8.68%  '
$ws.Range("D11").Value = $q + '1.839.05'
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("D12").Value = $q + '0.7273'
$ws.Range("E12").Value = '  -7.07%  '
$ws.Range("D13").Value = $q + '0.07193'
$ws.Range("E13").Value = '  -7.30%  '
$ws.Range("D14").Value = $q + '89.20'
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").Value = $q + '4.978'
$ws.Range("E15").Value = '  -1.06%  '
$ws.Range("D16").Value = $q + '1.001'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").Value = $q + '13.76'
$ws.Range("E17").Value = '  -1.34%  '
$ws.Range("D18").Value = $q + '1.000'
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").Value = $q + '0.000007896'
$ws.Range("E19").Value = '  -0.81%  '
$ws.Range("D20").Value = $q + '26.492.74'
$ws.Range("E20").Value = '  -0.29%  '
$ws.Range("D21").Value = $q + '2.082.74'
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("D22").Value = $q + '4.587'
$ws.Range("D23").Value = $q + '5.987'
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").Value = $q + '9.169'
$ws.Range("E24").Value = '  -3.08%  '
$ws.Range("D25").Value = $q + '142.97'
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").Value = $q + '2.158'
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("D27").Value = $q + '1.710'
$ws.Range("E27").Value = '  +1.77%  '
$ws.Range("D28").Value = $q + '16.97'
$ws.Range("E28").Value = '  -0.49%  '
$ws.Range("D29").Value = $q + '110.97'
$ws.Range("E29").Value = '  -1.18%  '
$ws.Range("D30").Value = $q + '4.258'
$ws.Range("E30").Value = '  +1.72%  '
$ws.Range("D31").Value = $q + '0.08902'
$ws.Range("E31").Value = '  +2.02%  '
$ws.Range("D32").Value = $q + '4.032'
$ws.Range("E32").Value = '  -1.75%  '
$ws.Range("D33").Value = $q + '0.04837'
$ws.Range("E33").Value = '  -0.69%  '
$ws.Range("D34").Value = $q + '2.903'
$ws.Range("E34").Value = '  +0.94%  '
$ws.Range("D35").Value = $q + '0.7255'
$ws.Range("E35").Value = '  +0.86%  '
$ws.Range("D36").Value = $q + '1.131'
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("D37").Value = $q + '3.095'
$ws.Range("E37").Value = '  -0.31%  '
$ws.Range("D38").Value = $q + '2.287'
$ws.Range("E38").Value = '  +0.26%  '
$ws.Range("D39").Value = $q + '0.01708'
$ws.Range("E39").Value = '  -4.34%  '
$ws.Range("E40").Value = '  -4.01%  '
$ws.Range("D41").Value = $q + '0.9018'
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").Value = $q + '107.35'
$ws.Range("E42").Value = '  -3.42%  '
$ws.Range("D43").Value = $q + '5.876'
$ws.Range("E43").Value = '  -1.03%  '
$ws.Range("D44").Value = $q + '1.000'
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").Value = $q + '7.416'
$ws.Range("E45").Value = '  -3.76%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = $q + '0.1243'
$ws.Range("E46").Value = '  +0.76%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = $q + '8.999'
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("D48").Value = $q + '0.4051'
$ws.Range("E48").Value = '  -2.97%  '
$ws.Range("D49").Value = $q + '34.79'
$ws.Range("E49").Value = '  -1.01%  '
$ws.Range("D50").Value = $q + '0.8943'
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = $q + '0.05736'
$ws.Range("E51").Value = '  -2.29%  '
